$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (current "数字" column), so that
# "数字" header splits into two headers: "数字Int" and "数字Long".
$ws.Range("B1").EntireColumn.Insert()

# Copy style (yellow fill) from A1 to the newly inserted B1 cell.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text.
$ws.Range("B1").Value = "数字Int"
$ws.Range("C1").Value = "数字Long"

$ws.Range("C1").Select()
